$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old hyperlink on D2 (deepak.n@ab-inbev.com mailto link) and
# reset that cell's formatting away from the "Hyperlink" style/font.
$ws.Hyperlinks.Delete()
$ws.Range("D2").Style = "Normal"

# Drop the now-unused "Hyperlink" cell style definition.
$wb.Styles.Item("Hyperlink").Delete()

# Update header row (row 1) and data row (row 2) to the new set of
# columns / test data (order chosen to match shared-string insertion order).
$ws.Range("B1").Value = "CreditCardNumber"
$ws.Range("D1").Value = "CVV"
$ws.Range("C1").Value = "ExpiryDate"
$ws.Range("C2").Value = "03/20"
$ws.Range("B2").Value = "4811 1111 1111 1114"
$ws.Range("E1").Value = "OTP"

$ws.Range("D2").Value = 123
$ws.Range("E2").Value = 112233

# Store credit card / expiry / CVV / OTP values as text so they keep
# leading context / formatting (e.g. "03/20") instead of being parsed.
$ws.Range("B2:E2").NumberFormat = "@"

# Widen column B so the longer credit card number is fully visible.
$ws.Columns("B").ColumnWidth = 15.8

# Move the active selection to B2 (was B10).
$ws.Range("B2").Select() | Out-Null
